$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Codelijsten")
$ws.Range("A1").Value = "test"
